# Add a hyperlink example to cell A8, similar to the other formatting
# examples already present on the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("A8")

# Create the hyperlink (this also applies the built-in "Hyperlink" cell
# style -- underline + theme color -- and registers the external
# relationship for xl/worksheets/_rels/sheet1.xml.rels).
$ws.Hyperlinks.Add($cell, "http://yahoo.com", "", "", "http://yahoo.com") | Out-Null

# The link text shown in the cell should be more descriptive than the bare
# address; set it after creating the hyperlink so the "display" attribute
# on the <hyperlink> element keeps referring to the address while the cell
# itself shows the friendly text.
$cell.Value = "http://yahoo.com Yahoo!"

# Move the active selection down one row, as happened after typing the
# new entry.
$ws.Range("A9").Select()
